$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (LangGraph): update Deadline
$ws.Range("F3").Value = "2025-09-14 18:42"

# Row 4 (FastAPI): rename task, update Deadline
$ws.Range("D4").Value = "Learn FastAPI"
$ws.Range("F4").Value = "2025-09-14 18:42"

# Row 5 (Subprocess): rename task, update Deadline
$ws.Range("D5").Value = "Learn Subprocess"
$ws.Range("F5").Value = "2025-09-14 18:42"
